$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Bug #59 fix: remove "somme_due"/"somme_payee" rows from the A column
# mini-table (subscriptions). "somme_due" (A6) is replaced by the former
# A8 value "message_abn", and the old A7/A8 cells are cleared so the data
# shifts up (table1 now spans A1:A6 instead of A1:A8).
$ws.Range("A6").Value = "message_abn"
$ws.Range("A7").Clear() | Out-Null
$ws.Range("A8").Clear() | Out-Null

# Resize the "subscriptions" table to match the shrunk data range.
$ws.ListObjects.Item("Table1").Resize($ws.Range("A1:A6"))

# --- Add the new "statut_act" legend block in column I (rows 19-23)
$ws.Range("I19").Value = "statut_act"
$ws.Range("I20").Value = "0=en cours"
$ws.Range("I21").Value = "1=validé"
$ws.Range("I22").Value = "2=payé"
$ws.Range("I23").Value = "3=refusé"

# --- Update the sheet view: move the active selection to A10.
$ws.Activate()
$ws.Range("A10").Select() | Out-Null

